$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ofmethod default value: "farneback" -> "mog2"
$ws.Range("B2").Value = "mog2"
$ws.Range("C2").Value = "mog2"

# Extend the comment on A2 describing the ofmethod options with the new
# mog / mog2 background-subtractor methods
$comment = $ws.Range("A2").Comment
$newCommentText = "hs: Horn-Schunck (opencv2 only)`nfarneback: Farneback method`nmog: background subtractor (opencv2 only)`nmog2: background subtractor"
$null = $comment.Text($newCommentText)

# Add new parameter rows for the GMM (mixture of gaussians) background
# subtractor (mog2) method
$ws.Range("A26").Value = "nhistory"
$ws.Range("B26").Value = 100
$ws.Range("C26").Value = 100

$ws.Range("A27").Value = "nmixtures"
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = 5

$ws.Range("A28").Value = "varThreshold"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1

# Move / persist the active selection to A30, matching the saved view state
$null = $ws.Range("A30").Select()

# Widen the sheet-tab area ratio in the saved window view
$win = $excel.ActiveWindow
$win.TabRatio = 0.983
